# Apply the "Adjusted for 2 times, ownership and lifetime." edit.
#
# Summary of the structural change (derived from the OOXML diff):
#   - Insert a new column A ("Reference") before the existing data,
#     shifting every existing column one to the right (A..V -> B..W).
#   - Insert a new column ("t_life_orig") between the shifted
#     "E_emb_orig" (now V) and "E_emb_star" (now W->X), i.e. insert
#     before the new column W, pushing "E_emb_star" from W to X.
#   - Append a new column Y ("t_life_star") after the new "E_emb_star".
#   - Rename the (now shifted) "t_orig"/"t_star" headers to
#     "t_own_orig"/"t_own_star".
#   - Fill in the new "Reference" column with "None yet" and the two
#     new "t_life_*" columns with ownership/lifetime year values.
#   - All existing cell comments are re-anchored one column to the
#     right (two for the last pair, because of the extra inserted
#     column), and four new comments ("years") are added for the new
#     t_life_orig / t_life_star cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Capture the existing comments (text + target new address) before
#    we disturb anything, then remove them. Comments stay bound to
#    their original row/column through cell moves/inserts in this
#    host, so we must explicitly recreate them at the shifted address.
# ---------------------------------------------------------------------
$oldToNew = @{
  "D2" = "E2";  "E2" = "F2";  "G2" = "H2";  "H2" = "I2";  "I2" = "J2"
  "M2" = "N2";  "N2" = "O2";  "O2" = "P2";  "P2" = "Q2";  "Q2" = "R2"
  "R2" = "S2";  "S2" = "T2";  "T2" = "U2";  "U2" = "V2";  "V2" = "X2"
  "D3" = "E3";  "E3" = "F3";  "G3" = "H3";  "H3" = "I3";  "I3" = "J3"
  "M3" = "N3";  "N3" = "O3";  "O3" = "P3";  "P3" = "Q3";  "Q3" = "R3"
  "R3" = "S3";  "S3" = "T3";  "T3" = "U3";  "U3" = "V3";  "V3" = "X3"
}

$pending = @()
foreach ($oldRef in $oldToNew.Keys) {
  $cmt = $ws.Range($oldRef).Comment
  if ($cmt) {
    $txt = $cmt.Text()
    $pending += [PSCustomObject]@{ NewRef = $oldToNew[$oldRef]; Text = $txt }
    $cmt.Delete()
  }
}

# ---------------------------------------------------------------------
# 2) Structural column inserts (this shifts cell values/styles/col
#    widths correctly in this host, even though it does not move the
#    comments - handled separately above/below).
# ---------------------------------------------------------------------
$ws.Columns("A:A").Insert(-4161)   # xlShiftToRight: insert new column A, old A..V -> B..W
$ws.Columns("W:W").Insert(-4161)   # insert new column W (old V/"E_emb_star" W -> X)

# ---------------------------------------------------------------------
# 3) Rename the shifted t_orig / t_star headers.
# ---------------------------------------------------------------------
$ws.Range("Q1").Value2 = "t_own_orig"
$ws.Range("S1").Value2 = "t_own_star"

# ---------------------------------------------------------------------
# 4) New "t_life_orig"/"t_life_star" columns (W, Y) and the new
#    "Reference" column (A).
# ---------------------------------------------------------------------
$ws.Range("W1").Value2 = "t_life_orig"
$ws.Range("Y1").Value2 = "t_life_star"

$ws.Range("A1").Value2 = "Reference"
$ws.Range("A2").Value2 = "None yet"
$ws.Range("A3").Value2 = "None yet"

$ws.Range("W2").Value2 = 14
$ws.Range("Y2").Value2 = 14

$ws.Range("W3").Value2 = 1.8
$ws.Range("Y3").Value2 = 10

# ---------------------------------------------------------------------
# 5) Re-create the shifted comments at their new address.
# ---------------------------------------------------------------------
foreach ($p in $pending) {
  $ws.Range($p.NewRef).AddComment($p.Text) | Out-Null
}

# ---------------------------------------------------------------------
# 6) Brand-new comments for the new ownership/lifetime columns.
# ---------------------------------------------------------------------
$yearsText = "Matthew Heun:" + [char]10 + "years"
$ws.Range("W2").AddComment($yearsText) | Out-Null
$ws.Range("Y2").AddComment($yearsText) | Out-Null
$ws.Range("W3").AddComment($yearsText) | Out-Null
$ws.Range("Y3").AddComment($yearsText) | Out-Null

# ---------------------------------------------------------------------
# 7) Misc view cleanup matching the saved file (selection moved to A4,
#    no frozen "topLeftCell").
# ---------------------------------------------------------------------
$ws.Range("A4").Select() | Out-Null
